# Apply updated KPI computations to the Productdata, Capacity and
# ProcessingTime sheets.

$wb = $excel.ActiveWorkbook

# --- Productdata sheet: update columns C (count) and E (computed KPI) ---
$ws = $wb.Worksheets.Item("Productdata")

$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 0.2222222222222222

$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 0.05141666666666665

$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 0.05138888888888889

$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 0.05130555555555555

$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 0.1026111111111111

$ws.Range("C7").Value = 5
$ws.Range("E7").Value = 0.05141666666666665

$ws.Range("C8").Value = 5
$ws.Range("E8").Value = 0.05138888888888889

$ws.Range("C9").Value = 5
$ws.Range("E9").Value = 0.1026111111111111

$ws.Range("C10").Value = 5
$ws.Range("E10").Value = 0.05130555555555555

$ws.Range("C11").Value = 0
$ws.Range("E11").Value = 0.04474999999999999

$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 0.04472222222222221

$ws.Range("C13").Value = 0
$ws.Range("E13").Value = 0.04463888888888888

$ws.Range("C14").Value = 0
$ws.Range("E14").Value = 0.08927777777777776

# --- Capacity sheet: update column B ---
$ws = $wb.Worksheets.Item("Capacity")

$ws.Range("B2").Value = 150
$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 40
$ws.Range("B5").Value = 50
$ws.Range("B6").Value = 20
$ws.Range("B7").Value = 10
$ws.Range("B8").Value = 40
$ws.Range("B9").Value = 40
$ws.Range("B10").Value = 20
$ws.Range("B11").Value = 20
$ws.Range("B12").Value = 50
$ws.Range("B13").Value = 10
$ws.Range("B14").Value = 40

# --- ProcessingTime sheet: update the diagonal cells ---
$ws = $wb.Worksheets.Item("ProcessingTime")

$ws.Range("B2").Value = 3
$ws.Range("C3").Value = 1
$ws.Range("D4").Value = 4
$ws.Range("E5").Value = 5
$ws.Range("F6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H8").Value = 4
$ws.Range("I9").Value = 4
$ws.Range("J10").Value = 2
$ws.Range("K11").Value = 2
$ws.Range("N14").Value = 2
